# journal_de_travail.xlsx — add a new "Utils" sheet holding the dropdown
# list that used to live in column H of the journal, point the data
# validation at it, log a new work session (2024-07-02) and tidy up.

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet -------------------------------------------
$journal = $wb.Worksheets.Item(1)
$journal.Name = "Journal"

# --- add the "Utils" sheet right after "Journal" and seed the list -------
$utils = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $journal)
$utils.Name = "Utils"

$utils.Range("B2").Value = "Analyse et état de l'art"
$utils.Range("B3").Value = "Réalisation du modèle"
$utils.Range("B4").Value = "Réalisation de l'application "
$utils.Range("B5").Value = "Tests et validations"
$utils.Range("B6").Value = "Gestion du projet, documentation et présentation"
$utils.Columns.Item(2).ColumnWidth = 45.43

# --- log the new working session on the Journal sheet --------------------
$journal.Range("A7").Value = 45475
$journal.Range("B7").Value = 0.375
$journal.Range("C7").Value = 0.63611111111111118

# --- the old helper list in column H is no longer needed -----------------
$journal.Range("H4:H8").ClearContents()
$journal.Columns.Item(8).AutoFit()

# --- repoint the dropdown validation at the new Utils sheet ---------------
$journal.Range("E1:E1048576").Validation.Delete()
$journal.Range("E1:E1048576").Validation.Add(3, 1, 1, "=Utils!`$B`$2:`$B`$6")

# --- restore cursor positions as last saved -------------------------------
$journal.Range("E12").Select()
$utils.Range("C17").Select()
$journal.Activate()
